$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.528.56"
$ws.Range("D3").Value = "1.750.86"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4473"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.05%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.094"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.96%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").Value = "1.749.42"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06378"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.859"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").Value = "27.571.78"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.082"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "1.948.80"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.090"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.084"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.663"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09049"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.553"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02298"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06021"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6369"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2086"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.945"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.384"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.776"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.724"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5903"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.149"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06859"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.56%  "
